$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for Wins, Losses, Ties using the same style as existing headers (AA1/AC1 etc, style index 1)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy style from an existing header cell (AC1) to the new header cells
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Fill in the win/loss/tie record for each data row (2 through 53)
for ($r = 2; $r -le 53; $r++) {
    $ws.Cells.Item($r, 30).Value = 97   # AD
    $ws.Cells.Item($r, 31).Value = 65   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
